# "push to git by dinh"
#
# The workbook started life as a little "stt / hoten / ngaysinh / tuoi"
# table on Sheet1, plus a second header+row table on Sheet2 and an empty
# Sheet3. The commit strips Sheet1 down to a tiny 3-cell "month" snippet
# (a running number, a "Tháng: " label, and a lone numeric value) and
# removes all the ad-hoc header styling (purple font / green fill /
# borders) that used to live on row 1. Sheet2 and Sheet3 are untouched.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: drop the two middle data rows -------------------------------
# Original column A was: stt-header, 3, 4, 5, 6 (rows 1-5).
# Row 3 (value 4) and what becomes the new row 4 (original row 5's value 6)
# are deleted, leaving the row that held "5" as the new row 3.
$ws1.Rows(3).Delete()
$ws1.Rows(4).Delete()

# --- Sheet1: clear the old header/table formatting ------------------------
# Wipes A1's purple-font/"stt" header plus the now write B1:C2 table cells
# (hoten/ngaysinh header + the row-2 numbers) entirely - values, shared
# string refs and the old purple-font/green-fill/border styling all go.
$ws1.Range("A1:C2").Clear()

# --- Sheet1: write the new 3-row content -----------------------------------
$ws1.Range("A1").Value = 1
$ws1.Range("A2").Value = "Tháng: "
# A3 already holds 5 (shifted up by the row deletes above) - matches target.
